$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 795, pushing the current
# rows 795:861 down to 797:863 (matches dimension growing from T861 to T863).
$ws.Range("A795:A796").EntireRow.Insert()

# --- New row 795: "1a plateado" entry dated 2023-06-29 (serial 45106) ---
$ws.Cells.Item(795, 1).Value = 4
$ws.Cells.Item(795, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(795, 3).Value = "Los Lagos"
$ws.Cells.Item(795, 4).Value = 45106
$ws.Cells.Item(795, 5).Value = 10
$ws.Cells.Item(795, 6).Value = "Fruta"
$ws.Cells.Item(795, 7).Value = 100102
$ws.Cells.Item(795, 8).Value = "Cítricos"
$ws.Cells.Item(795, 9).Value = 100102003
$ws.Cells.Item(795, 10).Value = "Limón"
$ws.Cells.Item(795, 11).Value = "Sin especificar"
$ws.Cells.Item(795, 12).Value = "1a plateado"
$ws.Cells.Item(795, 13).Value = 1000
$ws.Cells.Item(795, 14).Value = 14000
$ws.Cells.Item(795, 15).Value = 15000
$ws.Cells.Item(795, 16).Value = 14500
$ws.Cells.Item(795, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(795, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(795, 19).Value = 806
$ws.Cells.Item(795, 20).Value = 18

# --- New row 796: "2a plateado" entry dated 2023-06-29 (serial 45106) ---
$ws.Cells.Item(796, 1).Value = 4
$ws.Cells.Item(796, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(796, 3).Value = "Los Lagos"
$ws.Cells.Item(796, 4).Value = 45106
$ws.Cells.Item(796, 5).Value = 10
$ws.Cells.Item(796, 6).Value = "Fruta"
$ws.Cells.Item(796, 7).Value = 100102
$ws.Cells.Item(796, 8).Value = "Cítricos"
$ws.Cells.Item(796, 9).Value = 100102003
$ws.Cells.Item(796, 10).Value = "Limón"
$ws.Cells.Item(796, 11).Value = "Sin especificar"
$ws.Cells.Item(796, 12).Value = "2a plateado"
$ws.Cells.Item(796, 13).Value = 500
$ws.Cells.Item(796, 14).Value = 12000
$ws.Cells.Item(796, 15).Value = 12000
$ws.Cells.Item(796, 16).Value = 12000
$ws.Cells.Item(796, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(796, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(796, 19).Value = 667
$ws.Cells.Item(796, 20).Value = 18
